# Update need_to_buy.xlsx data with refreshed values from R (rows 2-15)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: row index (sheet row) => A,B,C,D,E,F
$data = @(
    @(2,  46006, 11113.0156395031, 10436.6553349417, 11232.26, 6978.20375910153, 257.608295585136),
    @(3,  46007, 11101.3582522158, 10222.1033172001, 11232.26, 7328.09507999179, 263.247433216329),
    @(4,  46008, 10790.9458229323, 9689.73833144721, 11232.26, 7086.58682267669, 231.002714755163),
    @(5,  46009, 10649.5736974066, 9870.31196130763, 11232.26, 6965.84627263553, 233.495759747632),
    @(6,  46010, 10697.3442480687, 9249.8214033901,  11232.26, 6998.49349955907, 209.002287622882),
    @(7,  46011, 8560.57524882407, 8717.30886156641, 11232.26, 7650.34018821925, 213.974543741069),
    @(8,  46012, 8469.04674334209, 8557.56092959685, 11232.26, 7641.94217211311, 206.968462571248),
    @(9,  46013, 9566.18303866457, 8948.77601078289, 11232.26, 8113.15183942236, 242.902827091885),
    @(10, 46014, 9566.18303866457, 9204.55257025269, 11232.26, 8113.15183942236, 253.56018373646),
    @(11, 46015, 9566.18303866457, 9392.45052026764, 11232.26, 8113.15183942236, 261.389264987083),
    @(12, 46016, 8471.80799626025, 9049.65841608714, 11232.26, 7768.77536493624, 232.757240875974),
    @(13, 46017, 8471.80799626025, 8472.62130634191, 11232.26, 7768.77536493624, 208.714027969923),
    @(14, 46018, 8547.24893097507, 9733.5221584338,  11232.26, 7775.95430298689, 261.550685892529),
    @(15, 46019, 8471.80799626025, 10014.4421063493, 11232.26, 7768.77536493624, 272.956561303562)
)

foreach ($r in $data) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
}
